$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "I" label that used to sit in A2
$ws.Range("A2").ClearContents()

# Remove the stray value in J8 / column J (solver no longer writes there)
$ws.Columns.Item(10).Delete()

# Move the active selection to K2, matching where the user clicked next
$ws.Range("K2").Select()
